$wb = $excel.ActiveWorkbook

# Add a new worksheet after the last existing sheet and name it "ghkn"
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "ghkn"

# Header row
$newSheet.Range("A1").Value = "Number"
$newSheet.Range("B1").Value = "Student ID"
$newSheet.Range("C1").Value = "Location"
$newSheet.Range("D1").Value = "Log Date"
$newSheet.Range("E1").Value = "Log Time"

# Data row
$newSheet.Range("A2").Value = 1
$newSheet.Range("B2").Value = "dfvk"
$newSheet.Range("C2").Value = "ghkn"
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "2025-04-05"
$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "21:21:45"
